# New data row (row 60) describing the LIDL supermarket location, added
# at the end of the "locais" table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A60").Value = "Supermercado"
$ws.Range("B60").Value = "LIDL"
$ws.Range("C60").Value = 55.942533277833803
$ws.Range("D60").Value = -3.22305196000559
$ws.Range("E60").Value = "supermarket.png"
$ws.Range("F60").Value = "Dalry Rd, Edinburgh EH11 2EF, Reino Unido"
$ws.Range("G60").Value = "Horário: 08:00–22:00"

# Match the "categoria" column formatting used by the rest of the table
# (left-aligned text), same as cell A59.
$ws.Range("A60").HorizontalAlignment = -4131

# Scroll the view down to the newly added row and select the same cell
# that was active when the workbook was last saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 51
$win.ScrollColumn = 1
$ws.Range("E66").Select()
